$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9780805706977844
$ws.Range("B1").Value = 2.817479610443115
$ws.Range("C1").Value = 3.045762777328491
$ws.Range("D1").Value = 3.650051832199097
$ws.Range("E1").Value = 1.587705492973328
